{"js": "// Update the worksheet date heading and the 25 division-problem cells\n// (5 data rows x 5 cols inside a 20-row table that also has 3 blank\n// spacer rows after every data row) to the next day's values.\n\nconst body = context.document.body;\n\n// 1) Date heading paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text === \"2024-06-24 Monday\") {\n  dateParagraph.insertText(\"2024-06-25 Tuesday\", \"Replace\");\n}\n\n// 2) Division-problem table cells, addressed by (row, col) so duplicate\n//    cell text (e.g. \"794\u00f75=158, 4\" appears twice before the edit) is\n//    replaced independently rather than via a global text search.\nconst table = body.tables.items[0];\n\n// Row indices of the 5 data rows (every 4th row: data row then 3 blanks).\nconst dataRowIndices = [0, 4, 8, 12, 16];\n\n// New values, in reading order, 5 columns per data row.\nconst newValues = [\n  [\"870\u00f75=174, 0\", \"681\u00f79=75, 6\", \"731\u00f75=146, 1\", \"950\u00f72=475, 0\", \"805\u00f75=161, 0\"],\n  [\"794\u00f75=158, 4\", \"754\u00f75=150, 4\", \"221\u00f76=36, 5\", \"576\u00f74=144, 0\", \"909\u00f75=181, 4\"],\n  [\"501\u00f76=83, 3\", \"300\u00f73=100, 0\", \"528\u00f77=75, 3\", \"503\u00f78=62, 7\", \"994\u00f77=142, 0\"],\n  [\"813\u00f74=203, 1\", \"487\u00f72=243, 1\", \"494\u00f76=82, 2\", \"833\u00f79=92, 5\", \"347\u00f75=69, 2\"],\n  [\"687\u00f77=98, 1\", \"738\u00f72=369, 0\", \"962\u00f73=320, 2\", \"874\u00f77=124, 6\", \"696\u00f79=77, 3\"],\n];\n\nfor (let r = 0; r < dataRowIndices.length; r++) {\n  const rowIdx = dataRowIndices[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIdx, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and the 25 division-problem cells\n# (5 data rows x 5 cols inside a 20-row table that also has 3 blank\n# spacer rows after every data row) to the next day's values.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph (first paragraph in the body).\n$dateParagraph = $d.Paragraphs(1)\nif ($dateParagraph.Range.Text -like \"2024-06-24 Monday*\") {\n    $dateParagraph.Range.Text = \"2024-06-25 Tuesday\"\n}\n\n# 2) Division-problem table cells, addressed by (row, col) so duplicate\n#    cell text (e.g. \"794\u00f75=158, 4\" appears twice before the edit) is\n#    replaced independently rather than via a global find/replace.\n$t = $d.Tables(1)\n\n# 1-based COM row numbers of the 5 data rows (every 4th row: data row\n# then 3 blanks) -> 1, 5, 9, 13, 17.\n$dataRows = @(1, 5, 9, 13, 17)\n\n# New values, in reading order, 5 columns per data row.\n$newValues = @(\n    @(\"870\u00f75=174, 0\", \"681\u00f79=75, 6\", \"731\u00f75=146, 1\", \"950\u00f72=475, 0\", \"805\u00f75=161, 0\"),\n    @(\"794\u00f75=158, 4\", \"754\u00f75=150, 4\", \"221\u00f76=36, 5\", \"576\u00f74=144, 0\", \"909\u00f75=181, 4\"),\n    @(\"501\u00f76=83, 3\", \"300\u00f73=100, 0\", \"528\u00f77=75, 3\", \"503\u00f78=62, 7\", \"994\u00f77=142, 0\"),\n    @(\"813\u00f74=203, 1\", \"487\u00f72=243, 1\", \"494\u00f76=82, 2\", \"833\u00f79=92, 5\", \"347\u00f75=69, 2\"),\n    @(\"687\u00f77=98, 1\", \"738\u00f72=369, 0\", \"962\u00f73=320, 2\", \"874\u00f77=124, 6\", \"696\u00f79=77, 3\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $rowNum = $dataRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($rowNum, $c).Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
